$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell reference and its new value. All target cells hold
# plain text (prices formatted with "." as thousands separators and
# percent strings with padding spaces), so force text storage via
# NumberFormat "@" before writing, then restore the default "Normal"
# style so no stray formatting is introduced.
$updates = @(
    @{ Cell = "D2"; Value = "67.633.62" }
    @{ Cell = "E2"; Value = "  -1.50%  " }
    @{ Cell = "D3"; Value = "3.328.00" }
    @{ Cell = "E3"; Value = "  -1.31%  " }
    @{ Cell = "D4"; Value = "0.998" }
    @{ Cell = "E4"; Value = "  +0.01%  " }
    @{ Cell = "D5"; Value = "582.45" }
    @{ Cell = "E5"; Value = "  -1.95%  " }
    @{ Cell = "D6"; Value = "175.60" }
    @{ Cell = "E6"; Value = "  -5.53%  " }
    @{ Cell = "E7"; Value = "  +0.09%  " }
    @{ Cell = "E8"; Value = "  -1.44%  " }
    @{ Cell = "D9"; Value = "3.324.29" }
    @{ Cell = "E9"; Value = "  -0.95%  " }
    @{ Cell = "D10"; Value = "0.177" }
    @{ Cell = "E10"; Value = "  -2.86%  " }
    @{ Cell = "D11"; Value = "0.575" }
    @{ Cell = "E11"; Value = "  -1.79%  " }
    @{ Cell = "D12"; Value = "45.43" }
    @{ Cell = "E12"; Value = "  -3.70%  " }
    @{ Cell = "E13"; Value = "  -3.85%  " }
    @{ Cell = "D14"; Value = "658.63" }
    @{ Cell = "E14"; Value = "  +2.67%  " }
    @{ Cell = "D15"; Value = "3.866.18" }
    @{ Cell = "E15"; Value = "  -1.11%  " }
    @{ Cell = "D16"; Value = "8.39" }
    @{ Cell = "E16"; Value = "  -1.67%  " }
    @{ Cell = "D17"; Value = "67.729.66" }
    @{ Cell = "E17"; Value = "  -1.54%  " }
    @{ Cell = "E18"; Value = "  -1.03%  " }
    @{ Cell = "D19"; Value = "3.337.95" }
    @{ Cell = "E19"; Value = "  -1.12%  " }
    @{ Cell = "D20"; Value = "17.37" }
    @{ Cell = "E20"; Value = "  -2.96%  " }
    @{ Cell = "D21"; Value = "10.92" }
    @{ Cell = "E21"; Value = "  -1.32%  " }
    @{ Cell = "D22"; Value = "0.888" }
    @{ Cell = "E22"; Value = "  -2.47%  " }
    @{ Cell = "D23"; Value = "5.44" }
    @{ Cell = "E23"; Value = "  +7.04%  " }
    @{ Cell = "D24"; Value = "17.01" }
    @{ Cell = "E24"; Value = "  -5.23%  " }
    @{ Cell = "D25"; Value = "99.57" }
    @{ Cell = "E25"; Value = "  +0.55%  " }
    @{ Cell = "D26"; Value = "3.85" }
    @{ Cell = "E26"; Value = "  -6.07%  " }
    @{ Cell = "D27"; Value = "2.66" }
    @{ Cell = "E27"; Value = "  -6.31%  " }
    @{ Cell = "D28"; Value = "9.24" }
    @{ Cell = "E28"; Value = "  -5.79%  " }
    @{ Cell = "D29"; Value = "33.59" }
    @{ Cell = "E29"; Value = "  +1.93%  " }
    @{ Cell = "D30"; Value = "7.38" }
    @{ Cell = "E30"; Value = "  +8.57%  " }
    @{ Cell = "D31"; Value = "8.42" }
    @{ Cell = "E31"; Value = "  -2.89%  " }
    @{ Cell = "D32"; Value = "594.18" }
    @{ Cell = "E32"; Value = "  -2.81%  " }
    @{ Cell = "D33"; Value = "10.98" }
    @{ Cell = "E33"; Value = "  -0.83%  " }
    @{ Cell = "E34"; Value = "  -1.10%  " }
    @{ Cell = "B35"; Value = "Dai" }
    @{ Cell = "C35"; Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai" }
    @{ Cell = "D35"; Value = "0.998" }
    @{ Cell = "E35"; Value = "  -0.17%  " }
    @{ Cell = "B36"; Value = "Maker" }
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr" }
    @{ Cell = "D36"; Value = "3.705.78" }
    @{ Cell = "E36"; Value = "  -7.06%  " }
    @{ Cell = "D37"; Value = "56.64" }
    @{ Cell = "E37"; Value = "  +1.01%  " }
    @{ Cell = "D38"; Value = "3.33" }
    @{ Cell = "E38"; Value = "  -9.61%  " }
    @{ Cell = "E39"; Value = "  +0.07%  " }
    @{ Cell = "D40"; Value = "33.86" }
    @{ Cell = "E40"; Value = "  +0.93%  " }
    @{ Cell = "D41"; Value = "2.62" }
    @{ Cell = "E41"; Value = "  -5.11%  " }
    @{ Cell = "D42"; Value = "3.11" }
    @{ Cell = "E42"; Value = "  -6.00%  " }
    @{ Cell = "E43"; Value = "  -2.91%  " }
    @{ Cell = "D44"; Value = "0.0₃0665" }
    @{ Cell = "E44"; Value = "  -5.52%  " }
    @{ Cell = "D45"; Value = "3.27" }
    @{ Cell = "E45"; Value = "  -4.78%  " }
    @{ Cell = "E46"; Value = "  -3.59%  " }
    @{ Cell = "B47"; Value = "ThetaToken" }
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta" }
    @{ Cell = "D47"; Value = "2.59" }
    @{ Cell = "E47"; Value = "  +0.06%  " }
    @{ Cell = "B48"; Value = "Stellar" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" }
    @{ Cell = "D48"; Value = "0.128" }
    @{ Cell = "E48"; Value = "  -1.29%  " }
    @{ Cell = "E49"; Value = "  +0.06%  " }
    @{ Cell = "D50"; Value = "1.34" }
    @{ Cell = "E50"; Value = "  +0.14%  " }
    @{ Cell = "D51"; Value = "127.15" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
